$wb = $excel.ActiveWorkbook

# --- Sheet1 (2EXT02_Protein): add example content rows 2-7 ---
$ws1 = $wb.Worksheets.Item(1)

$colB = @("Arg-C","Trypsin","TrypChymo","no cleavage","unspecific cleavage","2-idobenzoate")
$colE = @("nucleic acid","peptide","razor peptide","peptidoform ion","protein","metabolite")
$colH = @("emulsion","solution","suspension","gaseous sample state","solid sample state","liquid sample state")
$colK = @("H&E slide staining","IHC slide staining")
$colN = @("Good's buffer substance")

for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws1.Range("B$(2+$i)").Value = $colB[$i]
}
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws1.Range("E$(2+$i)").Value = $colE[$i]
}
for ($i = 0; $i -lt $colH.Length; $i++) {
    $ws1.Range("H$(2+$i)").Value = $colH[$i]
}
for ($i = 0; $i -lt $colK.Length; $i++) {
    $ws1.Range("K$(2+$i)").Value = $colK[$i]
}
for ($i = 0; $i -lt $colN.Length; $i++) {
    $ws1.Range("N$(2+$i)").Value = $colN[$i]
}

# Keep the "Unit [pH]" number formatting consistent down the newly used rows
$ws1.Range("R3:R7").NumberFormat = '0.00\ "pH"'

# Resize/extend the annotation table to cover the new rows
$tbl1 = $ws1.ListObjects.Item(1)
$tbl1.Resize($ws1.Range("A1:AA7"))

# --- Sheet2 (SwateTemplateMetadata): bump template version 1.1.3 -> 1.1.4 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B3").Value = "'1.1.4"
